# Append the 2025-03-18 price row to each "Solar_Prices" sheet,
# keeping the Date/Price columns formatted as text, matching the
# existing sheet layout (A = date string, B = price string).

$wb = $excel.ActiveWorkbook

$newRow = @{
    "N-Dense"                   = "40"
    "N-Type"                    = "43"
    "N-type Wafer"               = "1.19"
    "Cell Topcon 183mm"          = "0.298"
    "Module Topcon 183mm"        = "0.1"
    "Silver Rear_side"           = "5,440"
    "Silver Busbar front-side"   = "8,145"
    "Silver finger front-side"   = "8,195"
    "USD_CNY"                    = "7.2456"
}

foreach ($sheetName in $newRow.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    $dateCell = $ws.Range("A17")
    $priceCell = $ws.Range("B17")

    # Force text storage so "2025-03-18" / numeric-looking strings
    # are not reinterpreted as a date serial or a number.
    $dateCell.NumberFormat = "@"
    $priceCell.NumberFormat = "@"

    $dateCell.Value = "2025-03-18"
    $priceCell.Value = $newRow[$sheetName]
}
